$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost (Q) / Nord (R) coordinate values to whole numbers for rows 4-7
$ws.Range("Q4").Value = 575755
$ws.Range("R4").Value = 6703742

$ws.Range("Q5").Value = 575759
$ws.Range("R5").Value = 6703742

$ws.Range("Q6").Value = 575827
$ws.Range("R6").Value = 6703782

$ws.Range("Q7").Value = 575783
$ws.Range("R7").Value = 6703744

# Remove the placeholder "00:00" start/end time values (rows 4-6 only;
# row 7 keeps its real time values)
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
